$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7212032675743103
$ws.Range("B1").Value = 1.011853933334351
$ws.Range("C1").Value = 1.355706810951233
$ws.Range("D1").Value = 4.402457237243652
$ws.Range("E1").Value = 2.389496803283691
